$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend formatting (style) for the new rows 18-31 down from row 17 (same A=currency / B-E=number style)
$ws.Range("A17:E17").Copy()
$ws.Range("A18:E31").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Column B width: was 7.2 chars, now matches column A (9.6 chars).
# COM ColumnWidth quantizes to the nearest renderable grid step in this engine;
# 8.75 is the closest settable value that maps to the target column-A width bucket.
$ws.Columns.Item(2).ColumnWidth = 8.75

$ws.Range("A2").Value2 = '$BTC'
$ws.Range("B2").Value2 = 21084.02
$ws.Range("C2").Value2 = 406175601662.3179
$ws.Range("D2").Value2 = 23813331260.59324
$ws.Range("E2").Value2 = 19264618
$ws.Range("A3").Value2 = '$ETH'
$ws.Range("B3").Value2 = 1560.73
$ws.Range("C3").Value2 = 190992789187.8394
$ws.Range("D3").Value2 = 8048775248.455356
$ws.Range("E3").Value2 = 122373866.2178
$ws.Range("A4").Value2 = '$USDT'
$ws.Range("B4").Value2 = 1
$ws.Range("C4").Value2 = 66490594505.92384
$ws.Range("D4").Value2 = 32377041275.90711
$ws.Range("E4").Value2 = 73141766321.23428
$ws.Range("A5").Value2 = '$BNB'
$ws.Range("B5").Value2 = 298.83
$ws.Range("C5").Value2 = 47801065445.75539
$ws.Range("D5").Value2 = 563691854.0179112
$ws.Range("E5").Value2 = 159979963.5904293
$ws.Range("A6").Value2 = '$USDC'
$ws.Range("B6").Value2 = 1
$ws.Range("C6").Value2 = 43835901598.5809
$ws.Range("D6").Value2 = 2901189759.720745
$ws.Range("E6").Value2 = 43841459704.57129
$ws.Range("A7").Value2 = '$XRP'
$ws.Range("B7").Value2 = 0.39
$ws.Range("C7").Value2 = 19666967859.20561
$ws.Range("D7").Value2 = 1227605238.703969
$ws.Range("E7").Value2 = 99989156648
$ws.Range("A8").Value2 = '$BUSD'
$ws.Range("B8").Value2 = 1
$ws.Range("C8").Value2 = 16247560099.42515
$ws.Range("D8").Value2 = 8225416843.490399
$ws.Range("E8").Value2 = 16242596010.61406
$ws.Range("A9").Value2 = '$ADA'
$ws.Range("B9").Value2 = 0.35
$ws.Range("C9").Value2 = 12057715982.64891
$ws.Range("D9").Value2 = 343469550.7680789
$ws.Range("E9").Value2 = 35303937971.934
$ws.Range("A10").Value2 = '$DOGE'
$ws.Range("B10").Value2 = 0.08
$ws.Range("C10").Value2 = 11130913882.48604
$ws.Range("D10").Value2 = 461106665.3322925
$ws.Range("E10").Value2 = 132670764299.8941
$ws.Range("A11").Value2 = '$MATIC'
$ws.Range("B11").Value2 = 1.01
$ws.Range("C11").Value2 = 8835122977.265144
$ws.Range("D11").Value2 = 607386215.487919
$ws.Range("E11").Value2 = 10000000000
$ws.Range("A12").Value2 = '$SOL'
$ws.Range("B12").Value2 = 23.36
$ws.Range("C12").Value2 = 8661677693.50992
$ws.Range("D12").Value2 = 1261440340.750323
$ws.Range("E12").Value2 = 538755344.6046014
$ws.Range("A13").Value2 = '$DOT'
$ws.Range("B13").Value2 = 6
$ws.Range("C13").Value2 = 6934936947.359878
$ws.Range("D13").Value2 = 308301009.1270776
$ws.Range("E13").Value2 = 1271962399.3702
$ws.Range("A14").Value2 = '$LTC'
$ws.Range("B14").Value2 = 85.7
$ws.Range("C14").Value2 = 6176940690.369475
$ws.Range("D14").Value2 = 596893939.7728382
$ws.Range("E14").Value2 = 84000000
$ws.Range("A15").Value2 = '$DAI'
$ws.Range("B15").Value2 = 1
$ws.Range("C15").Value2 = 5833569481.472278
$ws.Range("D15").Value2 = 179521199.5863547
$ws.Range("E15").Value2 = 5835995685.216533
$ws.Range("A16").Value2 = '$SHIB'
$ws.Range("B16").Value2 = 0
$ws.Range("C16").Value2 = 5754003682.647244
$ws.Range("D16").Value2 = 311403153.8771881
$ws.Range("E16").Value2 = 589735030408322.8
$ws.Range("A17").Value2 = '$TRX'
$ws.Range("B17").Value2 = 0.06
$ws.Range("C17").Value2 = 5700945370.103837
$ws.Range("D17").Value2 = 228364886.564765
$ws.Range("E17").Value2 = 91874796149.61217
$ws.Range("A18").Value2 = '$AVAX'
$ws.Range("B18").Value2 = 16.79
$ws.Range("C18").Value2 = 5230115101.594823
$ws.Range("D18").Value2 = 450337835.0047283
$ws.Range("E18").Value2 = 416988132.7412243
$ws.Range("A19").Value2 = '$UNI'
$ws.Range("B19").Value2 = 6.55
$ws.Range("C19").Value2 = 4989400929.190395
$ws.Range("D19").Value2 = 107955481.4764907
$ws.Range("E19").Value2 = 1000000000
$ws.Range("A20").Value2 = '$WBTC'
$ws.Range("B20").Value2 = 21062.28
$ws.Range("C20").Value2 = 3865929525.85234
$ws.Range("D20").Value2 = 110113201.8696926
$ws.Range("E20").Value2 = 183547.54135694
$ws.Range("A21").Value2 = '$ATOM'
$ws.Range("B21").Value2 = 12.45
$ws.Range("C21").Value2 = 3565281172.250839
$ws.Range("D21").Value2 = 201600933.9508532
$ws.Range("E21").Value2 = 0
$ws.Range("A22").Value2 = '$LINK'
$ws.Range("B22").Value2 = 6.71
$ws.Range("C22").Value2 = 3406832494.598842
$ws.Range("D22").Value2 = 358928028.7217537
$ws.Range("E22").Value2 = 1000000000
$ws.Range("A23").Value2 = '$LEO'
$ws.Range("B23").Value2 = 3.47
$ws.Range("C23").Value2 = 3314007895.554636
$ws.Range("D23").Value2 = 2358953.69817772
$ws.Range("E23").Value2 = 985239504
$ws.Range("A24").Value2 = '$XMR'
$ws.Range("B24").Value2 = 171.25
$ws.Range("C24").Value2 = 3121561109.13796
$ws.Range("D24").Value2 = 93090638.92516409
$ws.Range("E24").Value2 = 18228028.59174703
$ws.Range("A25").Value2 = '$ETC'
$ws.Range("B25").Value2 = 22.35
$ws.Range("C25").Value2 = 3104496982.653294
$ws.Range("D25").Value2 = 244667031.6846894
$ws.Range("E25").Value2 = 210700000
$ws.Range("A26").Value2 = '$TON'
$ws.Range("B26").Value2 = 2.32
$ws.Range("C26").Value2 = 2837885231.549485
$ws.Range("D26").Value2 = 33460868.23835777
$ws.Range("E26").Value2 = 5047558528
$ws.Range("A27").Value2 = '$BCH'
$ws.Range("B27").Value2 = 124.77
$ws.Range("C27").Value2 = 2406407891.330467
$ws.Range("D27").Value2 = 219252338.7187677
$ws.Range("E27").Value2 = 19286118.75
$ws.Range("A28").Value2 = '$XLM'
$ws.Range("B28").Value2 = 0.09
$ws.Range("C28").Value2 = 2265896764.437149
$ws.Range("D28").Value2 = 63037359.2597696
$ws.Range("E28").Value2 = 50001787494.60719
$ws.Range("A29").Value2 = '$CRO'
$ws.Range("B29").Value2 = 0.08
$ws.Range("C29").Value2 = 1992142691.202807
$ws.Range("D29").Value2 = 88615368.63830242
$ws.Range("E29").Value2 = 30263013692
$ws.Range("A30").Value2 = '$OKB'
$ws.Range("B30").Value2 = 31.44
$ws.Range("C30").Value2 = 1886453324.786516
$ws.Range("D30").Value2 = 22055353.68384389
$ws.Range("E30").Value2 = 300000000
$ws.Range("A31").Value2 = '$NEAR'
$ws.Range("B31").Value2 = 'Erreur'
$ws.Range("C31").Value2 = 'Erreur'
$ws.Range("D31").Value2 = 'Erreur'
$ws.Range("E31").Value2 = 'Erreur'
